$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as Text so numeric-looking price strings
# (e.g. "0.9991", "242.80") are preserved verbatim instead of being
# auto-converted into floating point numbers by Excel.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "30.157.59"
$ws.Range("E2").Value = "  +1.03%  "

# Row 3
$ws.Range("D3").Value = "1.892.38"
$ws.Range("E3").Value = "  +0.32%  "

# Row 4
$ws.Range("D4").Value = "0.9991"
$ws.Range("E4").Value = "  -0.13%  "

# Row 5
$ws.Range("D5").Value = "0.7464"
$ws.Range("E5").Value = "  -0.15%  "

# Row 6
$ws.Range("D6").Value = "242.80"
$ws.Range("E6").Value = "  +0.10%  "

# Row 7
$ws.Range("D7").Value = "0.9998"
$ws.Range("E7").Value = "  -0.08%  "

# Row 8
$ws.Range("D8").Value = "0.3173"
$ws.Range("E8").Value = "  +1.87%  "

# Row 9
$ws.Range("D9").Value = "0.07248"
$ws.Range("E9").Value = "  +1.76%  "

# Row 10
$ws.Range("D10").Value = "25.03"
$ws.Range("E10").Value = "  -1.54%  "

# Row 11
$ws.Range("D11").Value = "0.08365"
$ws.Range("E11").Value = "  -1.26%  "

# Row 12
$ws.Range("D12").Value = "0.7630"
$ws.Range("E12").Value = "  +0.49%  "

# Row 13
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "5.454"
$ws.Range("E13").Value = "  +1.80%  "

# Row 14
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.899.86"
$ws.Range("E14").Value = "  -0.39%  "

# Row 15
$ws.Range("D15").Value = "93.08"
$ws.Range("E15").Value = "  -0.23%  "

# Row 16
$ws.Range("D16").Value = "6.191"
$ws.Range("E16").Value = "  +0.80%  "

# Row 17
$ws.Range("D17").Value = "30.203.38"
$ws.Range("E17").Value = "  +1.01%  "

# Row 18
$ws.Range("D18").Value = "251.32"
$ws.Range("E18").Value = "  +3.40%  "

# Row 19
$ws.Range("D19").Value = "13.68"
$ws.Range("E19").Value = "  -0.11%  "

# Row 20
$ws.Range("D20").Value = "0.000007875"
$ws.Range("E20").Value = "  +1.08%  "

# Row 21
$ws.Range("D21").Value = "2.160.09"
$ws.Range("E21").Value = "  +0.08%  "

# Row 22
$ws.Range("D22").Value = "0.9995"
$ws.Range("E22").Value = "  +0.03%  "

# Row 23
$ws.Range("D23").Value = "8.035"
$ws.Range("E23").Value = "  +0.35%  "

# Row 24
$ws.Range("D24").Value = "0.9992"
$ws.Range("E24").Value = "  -0.11%  "

# Row 25
$ws.Range("D25").Value = "0.1592"
$ws.Range("E25").Value = "  +0.08%  "

# Row 26
$ws.Range("D26").Value = "9.322"
$ws.Range("E26").Value = "  -0.58%  "

# Row 27
$ws.Range("D27").Value = "164.24"
$ws.Range("E27").Value = "  +1.16%  "

# Row 28
$ws.Range("D28").Value = "18.82"
$ws.Range("E28").Value = "  +0.49%  "

# Row 29
$ws.Range("D29").Value = "2.075"
$ws.Range("E29").Value = "  +2.40%  "

# Row 30
$ws.Range("D30").Value = "1.477"
$ws.Range("E30").Value = "  -2.24%  "

# Row 31
$ws.Range("D31").Value = "4.611"
$ws.Range("E31").Value = "  +3.19%  "

# Row 32
$ws.Range("D32").Value = "1.538"
$ws.Range("E32").Value = "  +0.53%  "

# Row 33
$ws.Range("D33").Value = "4.232"
$ws.Range("E33").Value = "  +3.24%  "

# Row 34
$ws.Range("D34").Value = "0.05412"
$ws.Range("E34").Value = "  +0.32%  "

# Row 35
$ws.Range("D35").Value = "1.257"
$ws.Range("E35").Value = "  +1.61%  "

# Row 36
$ws.Range("D36").Value = "0.7684"
$ws.Range("E36").Value = "  +3.38%  "

# Row 37
$ws.Range("D37").Value = "0.9935"
$ws.Range("E37").Value = "  -1.02%  "

# Row 38
$ws.Range("D38").Value = "2.715"
$ws.Range("E38").Value = "  +0.17%  "

# Row 39
$ws.Range("D39").Value = "0.01976"
$ws.Range("E39").Value = "  +2.26%  "

# Row 40
$ws.Range("D40").Value = "2.772"
$ws.Range("E40").Value = "  +0.15%  "

# Row 41
$ws.Range("D41").Value = "0.4587"
$ws.Range("E41").Value = "  +3.03%  "

# Row 42
$ws.Range("D42").Value = "1.105.26"
$ws.Range("E42").Value = "  +1.57%  "

# Row 43
$ws.Range("D43").Value = "6.104"
$ws.Range("E43").Value = "  +0.53%  "

# Row 44
$ws.Range("D44").Value = "73.13"
$ws.Range("E44").Value = "  +0.68%  "

# Row 45
$ws.Range("D45").Value = "0.8705"
$ws.Range("E45").Value = "  +1.25%  "

# Row 46
$ws.Range("D46").Value = "104.50"
$ws.Range("E46").Value = "  +1.90%  "

# Row 47
$ws.Range("E47").Value = "  +0.02%  "

# Row 48
$ws.Range("D48").Value = "1.875"
$ws.Range("E48").Value = "  +0.77%  "

# Row 49
$ws.Range("D49").Value = "7.639"
$ws.Range("E49").Value = "  -0.42%  "

# Row 50
$ws.Range("D50").Value = "9.646"
$ws.Range("E50").Value = "  -0.77%  "

# Row 51
$ws.Range("D51").Value = "2.053.80"
$ws.Range("E51").Value = "  -0.14%  "

# Restore column D to its original (default/"Normal") style now that the
# text values have been written, so no stray style/number-format diff
# is left behind on unrelated cells.
$dRange.Style = "Normal"
